$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous data row (row 9) onto the new row 10
# so the new row's label cell (A10) picks up the same bold/border/center style.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 6.4
$ws.Range("C10").Value = 35.4
$ws.Range("D10").Value = 20.5
$ws.Range("E10").Value = 10.7
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 8
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 14.4
$ws.Range("K10").Value = 8.6
$ws.Range("L10").Value = 19.5
$ws.Range("M10").Value = 11.2
$ws.Range("N10").Value = 12.1
$ws.Range("O10").Value = 6.8
$ws.Range("P10").Value = 8.4
$ws.Range("Q10").Value = 8.300000000000001
$ws.Range("R10").Value = 3.5
